$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "COBWEB" right after the existing "COLUMBUS" row (row 6),
# pushing the remaining rows down by one.
$ws.Rows.Item(7).Insert()

# The final (target) data, in row order, for columns A:E (rows 2-28).
$data = @(
    ,@("AppDoctor", "No", "None", "Bytecode", 2014)
    ,@("ARES", "Yes", "Auto", "Bytecode,OS", 2018)
    ,@("Brahmastra", "Partly", "Manual", "Bytecode", 2014)
    ,@("CAR", "No", "Auto", "OS", 2022)
    ,@("COLUMBUS", "No", "None", "Bytecode", 2023)
    ,@("COBWEB", "No", "None", "Source code", 2019)
    ,@("ConDroid", "Yes", "Auto", "Bytecode", 2014)
    ,@("CrashFuzzer", "No", "Auto", "None", 2016)
    ,@("Crashscope", "No", "Auto", "Bytecode", 2016)
    ,@("DALT", "No", "Auto", "Bytecode", 2022)
    ,@("DirectDroid", "No", "Manual", "None", 2018)
    ,@("Droid-ANTIRM", "No", "Manual", "Bytecode", 2017)
    ,@("DroidFuzzer", "Yes", "Auto", "Unspecified", 2013)
    ,@("EHBDroid", "Yes", "None", "Frame", 2016)
    ,@("FAX", "Yes", "Auto", "Bytecode", 2020)
    ,@("FuzzDroid", "Yes", "Manual", "Bytecode", 2017)
    ,@("GroddDroid", "Yes", "Manual", "Bytecode", 2015)
    ,@("Harvester", "No", "Manual", "Bytecode", 2016)
    ,@("IntelliDroid", "Yes", "Manual", "OS", 2016)
    ,@("IntentFuzzer", "Yes", "Auto", "Framework", 2014)
    ,@("Malton", "No", "M/Auto", "None", 2017)
    ,@("MoSSOT", "No", "Auto", "None", 2019)
    ,@("Null intent fuzzer", "Yes", "Auto", "None", 2009)
    ,@("Sasnauskas et al.", "No", "Auto", "None", 2014)
    ,@("SIEVE", "No", "Auto", "Bytecode", 2020)
    ,@("SMARTGEN", "No", "Auto", "Bytecode", 2017)
    ,@("Snowdrop", "No", "Auto", "Bytecode", 2017)
)

# Write column-by-column (B, C, A, D, E) so that brand-new shared strings are
# interned in the same relative order as the authored workbook.
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $r = $r + 1
}
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# The trailing "Usman et al." row has shifted down to row 29 after the insert
# above; remove it entirely since it's been dropped from the sheet.
$ws.Rows.Item(29).Delete()

# Update the view: drop the frozen/scrolled topLeftCell and move the active
# selection to F12.
$ws.Range("F12").Select()
